# Pagination 객체 구현
# Update the API design sheet: the "게시글 수정" (update article) row's URL
# no longer includes the "{article.pk}" path parameter.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 8 corresponds to the PUT /board/free/article... "게시글 수정" entry.
# Its URL cell (C8) changes from "/board/free/article/{article.pk}"
# to "/board/free/article".
$ws.Range("C8").Value = "/board/free/article"

# Update the active selection to match the author's final cursor position.
$ws.Range("E13").Select()
